# Applies the "Adapted for mobile devices" commit:
#  - Updates existing lexeme row 114 (W_zmār / "root") Groups value from
#    "Prār " to "Prār;Kṣur"
#  - Appends 5 new lexeme rows (264-268) at the bottom of the Affixation sheet
#  - Scrolls the view down to the newly-added rows (mirrors the author's
#    saved scroll position / selection in the workbook)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 114, column G (Groups) ---
$ws.Cells.Item(114, 7).Value = 'Prār;Kṣur'

# --- Append new rows at the end of the table ---
$ws.Cells.Item(264, 1).Value = 'I_ynālek'
$ws.Cells.Item(264, 2).Value = 'elsewhere'
$ws.Cells.Item(264, 3).Value = 'R_ālăk'
$ws.Cells.Item(264, 4).Value = 'PREF'

$ws.Cells.Item(265, 1).Value = 'W_wrassi'
$ws.Cells.Item(265, 2).Value = 'suffering'
$ws.Cells.Item(265, 3).Value = 'V_v`ras'
$ws.Cells.Item(265, 4).Value = 'INF'

$ws.Cells.Item(266, 1).Value = 'R_waro'
$ws.Cells.Item(266, 2).Value = 'stinky'
$ws.Cells.Item(266, 3).Value = 'K_v`är'
$ws.Cells.Item(266, 4).Value = 'REL1'

$ws.Cells.Item(267, 1).Value = 'R_warom'
$ws.Cells.Item(267, 2).Value = 'stinky'
$ws.Cells.Item(267, 3).Value = 'M_war'
$ws.Cells.Item(267, 4).Value = 'REL1'

$ws.Cells.Item(268, 1).Value = 'M_walāṃ'
$ws.Cells.Item(268, 2).Value = 'tent?'
$ws.Cells.Item(268, 3).Value = 'V_vāl'
$ws.Cells.Item(268, 4).Value = 'INAN4'

# --- Scroll the visible window down to the newly added rows ---
$ws.Application.Goto($ws.Range("D268"), $true)
$ws.Range("D268").Select()
